$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-23 Sunday" "2025-11-24 Monday"

Replace-Text "58÷8=" "71÷5="
Replace-Text "16÷7=" "99÷2="
Replace-Text "50÷8=" "26÷8="
Replace-Text "53÷4=" "14÷3="
Replace-Text "62÷7=" "15÷8="
Replace-Text "45÷6=" "89÷3="
Replace-Text "78÷8=" "51÷6="
Replace-Text "20÷6=" "58÷9="
Replace-Text "14÷7=" "37÷8="
Replace-Text "90÷6=" "49÷7="
Replace-Text "40÷3=" "67÷5="
Replace-Text "33÷7=" "18÷6="
Replace-Text "41÷4=" "57÷7="
Replace-Text "24÷7=" "55÷3="
Replace-Text "71÷8=" "32÷7="
Replace-Text "72÷3=" "70÷5="
Replace-Text "11÷4=" "12÷3="
Replace-Text "78÷6=" "66÷8="
Replace-Text "62÷2=" "13÷6="
Replace-Text "93÷2=" "97÷9="
Replace-Text "73÷4=" "74÷7="
Replace-Text "90÷4=" "71÷5="
Replace-Text "33÷8=" "66÷9="
Replace-Text "13÷4=" "48÷4="
Replace-Text "92÷7=" "63÷5="
